# 2016-5.xlsx update:
#  - Append two new measurement rows to the "VRKT" sheet.
#  - Append four new measurement rows to the "TSRT" sheet.
#  - Move the "TSRT" tab from the 2nd position to the last position
#    (after KZBG, AGMS, TBL01).
#
# Numeric-looking readings (e.g. "15.90", "3.10", "92.00") must be stored
# as TEXT (shared strings), not coerced to numbers, so we build them as
# ="..." text formulas first and then Copy / PasteSpecial(xlPasteValues)
# to collapse them down to literal string cells while keeping the
# existing cell style (s="1") intact.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# VRKT: add rows 2-3
# ---------------------------------------------------------------
$vrkt = $wb.Worksheets.Item("VRKT")

$vrkt.Range("A2").Value = 29
$vrkt.Range("B2").Formula = "=""15.90"""
$vrkt.Range("C2").Formula = "=""3.10"""
$vrkt.Range("D2").Value = "*"
$vrkt.Range("E2").Value = "*"
$vrkt.Range("F2").Value = "*"
$vrkt.Range("G2").Value = "*"

$vrkt.Range("A3").Value = 30
$vrkt.Range("B3").Formula = "=""62.58"""
$vrkt.Range("C3").Formula = "=""4.88"""
$vrkt.Range("D3").Formula = "=""26.92"""
$vrkt.Range("E3").Formula = "=""33.05"""
$vrkt.Range("F3").Formula = "=""7.38"""
$vrkt.Range("G3").Formula = "=""0.41"""

$vrktRange = $vrkt.Range("A2:G3")
$vrktRange.Copy()
$vrktRange.PasteSpecial(-4163)

# ---------------------------------------------------------------
# TSRT: add rows 2-5
# ---------------------------------------------------------------
$tsrt = $wb.Worksheets.Item("TSRT")

$tsrt.Range("A2").Value = 25
$tsrt.Range("B2").Formula = "=""241.92"""
$tsrt.Range("C2").Formula = "=""243.42"""
$tsrt.Range("D2").Value = "*"
$tsrt.Range("E2").Formula = "=""385.36"""
$tsrt.Range("F2").Value = "*"
$tsrt.Range("G2").Formula = "=""32.58"""

$tsrt.Range("A3").Value = 26
$tsrt.Range("B3").Formula = "=""50.62"""
$tsrt.Range("C3").Formula = "=""110.21"""
$tsrt.Range("D3").Formula = "=""141.43"""
$tsrt.Range("E3").Formula = "=""129.32"""
$tsrt.Range("F3").Formula = "=""92.00"""
$tsrt.Range("G3").Formula = "=""0.09"""

$tsrt.Range("A4").Value = 27
$tsrt.Range("B4").Formula = "=""48.17"""
$tsrt.Range("C4").Formula = "=""104.12"""
$tsrt.Range("D4").Formula = "=""127.04"""
$tsrt.Range("E4").Formula = "=""115.67"""
$tsrt.Range("F4").Formula = "=""93.50"""
$tsrt.Range("G4").Formula = "=""0.09"""

$tsrt.Range("A5").Value = 28
$tsrt.Range("B5").Formula = "=""49.83"""
$tsrt.Range("C5").Formula = "=""107.13"""
$tsrt.Range("D5").Formula = "=""121.81"""
$tsrt.Range("E5").Formula = "=""118.97"""
$tsrt.Range("F5").Formula = "=""93.62"""
$tsrt.Range("G5").Formula = "=""0.07"""

$tsrtRange = $tsrt.Range("A2:G5")
$tsrtRange.Copy()
$tsrtRange.PasteSpecial(-4163)

# ---------------------------------------------------------------
# Reorder tabs: move TSRT to the end, after VRKT, KZBG, AGMS, TBL01.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tsrt.Move($null, $lastSheet)

$excel.CutCopyMode = $false
